{"js": "// Apply the \"Added many more features\" edits to the Lost Relics review.\n// Each change is a straightforward text replacement; we locate the old\n// text with Body.search (exact, case-sensitive match) and replace every\n// matching run via Range.insertText(..., Word.InsertLocation.replace).\n\nconst replacements = [\n  // Title / H1 heading (also reused verbatim later in a bold run).\n  {\n    find: \"Play Lost Relics for Free - Unique Cluster Payout System\",\n    replace: \"Play Lost Relics Free - A Unique and Engaging Slot Game\",\n  },\n  // \"What we like\" bullets.\n  {\n    find: \"Unique gameplay with 5x5 grid and cluster payout system\",\n    replace: \"Unique cluster payout system\",\n  },\n  {\n    find: \"Great attention to detail in graphics, symbols and protagonist design\",\n    replace: \"Excellent visuals and design\",\n  },\n  {\n    find: \"Historically accurate symbols that add to the overall experience\",\n    replace: \"Attention to detail in symbol design\",\n  },\n  {\n    find: \"Fresh take on historical-themed slots\",\n    replace: \"Engaging gameplay mechanics and special features\",\n  },\n  // \"What we don't like\" bullets.\n  {\n    find: \"No progressive jackpot\",\n    replace: \"Not available on all online casino platforms\",\n  },\n  {\n    find: \"May not appeal to those who prefer traditional slots\",\n    replace: \"May take some time to understand cluster payout system\",\n  },\n  // Closing italic summary line.\n  {\n    find:\n      \"Explore ancient cultures with NetEnt's Lost Relics. Play for free and discover its unique cluster payout system and exciting special features.\",\n    replace:\n      \"Play Lost Relics for free and explore ancient history with unique cluster payout system.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Lost Relics review.\n# Each change is a straightforward text replacement; Find/Replace (ReplaceAll)\n# over the whole document content handles both single and duplicate\n# occurrences (the title text appears twice - once in the H1 heading and\n# once in a bold run further down - and both need the same new text).\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replaceText\n  $find.Execute(\n    $findText,\n    $true,   # MatchCase\n    $false,  # MatchWholeWord\n    $false,  # MatchWildcards\n    $false,  # MatchSoundsLike\n    $false,  # MatchAllWordForms\n    $true,   # Forward\n    1,       # Wrap = wdFindContinue\n    $false,  # Format\n    $replaceText,\n    2        # Replace = wdReplaceAll\n  ) | Out-Null\n}\n\n# Title / H1 heading (also reused verbatim later in a bold run).\nReplace-AllText \"Play Lost Relics for Free - Unique Cluster Payout System\" \"Play Lost Relics Free - A Unique and Engaging Slot Game\"\n\n# \"What we like\" bullets.\nReplace-AllText \"Unique gameplay with 5x5 grid and cluster payout system\" \"Unique cluster payout system\"\nReplace-AllText \"Great attention to detail in graphics, symbols and protagonist design\" \"Excellent visuals and design\"\nReplace-AllText \"Historically accurate symbols that add to the overall experience\" \"Attention to detail in symbol design\"\nReplace-AllText \"Fresh take on historical-themed slots\" \"Engaging gameplay mechanics and special features\"\n\n# \"What we don't like\" bullets.\nReplace-AllText \"No progressive jackpot\" \"Not available on all online casino platforms\"\nReplace-AllText \"May not appeal to those who prefer traditional slots\" \"May take some time to understand cluster payout system\"\n\n# Closing italic summary line.\nReplace-AllText \"Explore ancient cultures with NetEnt's Lost Relics. Play for free and discover its unique cluster payout system and exciting special features.\" \"Play Lost Relics for free and explore ancient history with unique cluster payout system.\"\n"}
